# Added team record (Wins/Losses/Ties) to the data sheet, alongside the
# existing roster columns (A:AC), in new columns AD:AF.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (bold, bordered, centered) from A1
# onto the three new header cells before setting their text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the team's overall record: 73 wins, 89 losses, 0 ties.
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 73  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 89  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
